# Apply the "PO Forecast" update to the B0CKHM2ZQ6 PO data workbook.
#
# 1) Rename the "Requested quantity" header on the two existing sheets
#    ("Weekly Quantity" -> Weekly_PO_Qty, "Monthly Trend" -> Monthly_PO_Qty).
# 2) Add a new "PO Forecast" worksheet (as the last tab) with the forecast
#    table (ds / PO_Forecast / yhat_lower / yhat_upper), reusing the same
#    header/date formatting already used on the other sheets.

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- 1) Update existing header labels -------------------------------------
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet ------------------------------------
# Duplicate an existing sheet (rather than inserting a blank one) so the new
# tab inherits the same sheet-level setup (outline props, page margins, ...)
# used throughout this workbook, then wipe its contents before filling it
# with the forecast data.
$wsMonthly.Copy($null, $wsMonthly)
$wsForecast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast.Name = "PO Forecast"
$wsForecast.Cells.Clear()

$headers = @("ds", "PO_Forecast", "yhat_lower", "yhat_upper")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $wsForecast.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$data = @(
    @(45214.99999999999, 89,  36.01009206990511,  137.5491389524055),
    @(45228.99999999999, 97,  45.96813500303507,  150.7136481936353),
    @(45235.99999999999, 101, 48.67899831613723,  152.1499844686127),
    @(45242.99999999999, 105, 50.3236440772446,   157.0839615687336),
    @(45249.99999999999, 109, 57.50347786745576,  163.2758241723637),
    @(45298.99999999999, 138, 83.32567806601364,  187.179445053441),
    @(45305.99999999999, 142, 90.9261967006817,   191.9681068148047),
    @(45312.99999999999, 146, 92.3736753043548,   196.4288803594593),
    @(45319.99999999999, 150, 98.24764906396794,  200.9657352894827),
    @(45326.99999999999, 154, 100.2389857283445,  200.0135495367424),
    @(45333.99999999999, 158, 103.8605121117479,  207.624203785168),
    @(45340.99999999999, 162, 113.0709751330943,  213.8050016294895),
    @(45347.99999999999, 166, 114.3158145356619,  218.3037384027098),
    @(45354.99999999999, 170, 115.9130639873493,  226.8515136291911)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsForecast.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Reuse the header style (bold/centered) and the date-number-format style
# already defined in the workbook so the new sheet matches the others.
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A15").PasteSpecial(-4122) # xlPasteFormats
